$d = $word.ActiveDocument

$d.Content.Find.Execute("483×6=", $true, $false, $false, $false, $false, $true, 1, $false, "218×3=", 2) | Out-Null
$d.Content.Find.Execute("915×3=", $true, $false, $false, $false, $false, $true, 1, $false, "875×7=", 2) | Out-Null
$d.Content.Find.Execute("567×9=", $true, $false, $false, $false, $false, $true, 1, $false, "380×5=", 2) | Out-Null
$d.Content.Find.Execute("501×9=", $true, $false, $false, $false, $false, $true, 1, $false, "967×4=", 2) | Out-Null
$d.Content.Find.Execute("837×9=", $true, $false, $false, $false, $false, $true, 1, $false, "727×6=", 2) | Out-Null
$d.Content.Find.Execute("239×6=", $true, $false, $false, $false, $false, $true, 1, $false, "442×5=", 2) | Out-Null
$d.Content.Find.Execute("543×4=", $true, $false, $false, $false, $false, $true, 1, $false, "912×7=", 2) | Out-Null
$d.Content.Find.Execute("276×2=", $true, $false, $false, $false, $false, $true, 1, $false, "323×5=", 2) | Out-Null
$d.Content.Find.Execute("794×9=", $true, $false, $false, $false, $false, $true, 1, $false, "639×7=", 2) | Out-Null
$d.Content.Find.Execute("635×2=", $true, $false, $false, $false, $false, $true, 1, $false, "586×3=", 2) | Out-Null
$d.Content.Find.Execute("513×5=", $true, $false, $false, $false, $false, $true, 1, $false, "998×4=", 2) | Out-Null
$d.Content.Find.Execute("984×2=", $true, $false, $false, $false, $false, $true, 1, $false, "411×4=", 2) | Out-Null
$d.Content.Find.Execute("370×5=", $true, $false, $false, $false, $false, $true, 1, $false, "328×9=", 2) | Out-Null
$d.Content.Find.Execute("797×4=", $true, $false, $false, $false, $false, $true, 1, $false, "482×5=", 2) | Out-Null
$d.Content.Find.Execute("975×7=", $true, $false, $false, $false, $false, $true, 1, $false, "881×5=", 2) | Out-Null
$d.Content.Find.Execute("143×2=", $true, $false, $false, $false, $false, $true, 1, $false, "669×2=", 2) | Out-Null
$d.Content.Find.Execute("691×9=", $true, $false, $false, $false, $false, $true, 1, $false, "249×6=", 2) | Out-Null
$d.Content.Find.Execute("190×5=", $true, $false, $false, $false, $false, $true, 1, $false, "742×2=", 2) | Out-Null
$d.Content.Find.Execute("920×9=", $true, $false, $false, $false, $false, $true, 1, $false, "913×5=", 2) | Out-Null
$d.Content.Find.Execute("839×4=", $true, $false, $false, $false, $false, $true, 1, $false, "321×5=", 2) | Out-Null
$d.Content.Find.Execute("191×2=", $true, $false, $false, $false, $false, $true, 1, $false, "491×5=", 2) | Out-Null
$d.Content.Find.Execute("139×6=", $true, $false, $false, $false, $false, $true, 1, $false, "370×2=", 2) | Out-Null
$d.Content.Find.Execute("981×3=", $true, $false, $false, $false, $false, $true, 1, $false, "346×3=", 2) | Out-Null
$d.Content.Find.Execute("283×2=", $true, $false, $false, $false, $false, $true, 1, $false, "145×3=", 2) | Out-Null
$d.Content.Find.Execute("696×5=", $true, $false, $false, $false, $false, $true, 1, $false, "414×3=", 2) | Out-Null
